$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("AdmSR")

$updates = @(
    @{Row=3; A=20150926; E=13},
    @{Row=4; A=20150927; E=14},
    @{Row=5; A=20150928; E=18},
    @{Row=6; A=20150929; E=20},
    @{Row=7; A=20150930; E=10},
    @{Row=8; A=20150931; E=16},
    @{Row=9; A=20150932; E=$null},
    @{Row=10; A=20150933; E=10},
    @{Row=11; A=20150934; E=18},
    @{Row=12; A=20150935; E=14},
    @{Row=13; A=20150936; E=15},
    @{Row=14; A=20150937; E=12},
    @{Row=15; A=20150938; E=18},
    @{Row=16; A=20150939; E=19},
    @{Row=17; A=20150940; E=5},
    @{Row=18; A=20150941; E=11},
    @{Row=19; A=20150942; E=5},
    @{Row=20; A=20150943; E=14},
    @{Row=21; A=20150944; E=19},
    @{Row=22; A=20150945; E=8},
    @{Row=23; A=20150946; E=12},
    @{Row=24; A=20150947; E=14},
    @{Row=25; A=20150948; E=$null},
    @{Row=26; A=20150949; E=10},
    @{Row=27; A=20150950; E=14},
    @{Row=28; A=20150951; E=6},
    @{Row=29; A=20150952; E=7},
    @{Row=30; A=20150953; E=9},
    @{Row=31; A=20150954; E=8},
    @{Row=32; A=20150955; E=$null},
    @{Row=33; A=20150956; E=9},
    @{Row=34; A=20150957; E=18},
    @{Row=35; A=20150958; E=7},
    @{Row=36; A=20150959; E=$null},
    @{Row=37; A=20150960; E=14},
    @{Row=38; A=20150961; E=14},
    @{Row=39; A=20150962; E=20},
    @{Row=40; A=20150963; E=$null},
    @{Row=41; A=20150964; E=20},
    @{Row=42; A=20150965; E=8},
    @{Row=43; A=20150966; E=5},
    @{Row=44; A=20150967; E=10},
    @{Row=45; A=20150968; E=20},
    @{Row=46; A=20150969; E=$null},
    @{Row=47; A=20150970; E=7},
    @{Row=48; A=20150971; E=18},
    @{Row=49; A=20150972; E=18},
    @{Row=50; A=20150973; E=5},
    @{Row=51; A=20150974; E=9},
    @{Row=52; A=20150975; E=12},
    @{Row=53; A=20150976; E=13},
    @{Row=54; A=20150977; E=15},
    @{Row=55; A=20150978; E=12},
    @{Row=56; A=20150979; E=8},
    @{Row=57; A=20150980; E=7},
    @{Row=58; A=20150981; E=15},
    @{Row=59; A=20150982; E=8},
    @{Row=60; A=20150983; E=5},
    @{Row=61; A=20150984; E=16},
    @{Row=62; A=20150985; E=6},
    @{Row=63; A=20150986; E=18}
)

foreach ($u in $updates) {
    $ws.Cells.Item($u.Row, 1).Value = $u.A
    if ($null -ne $u.E) {
        $ws.Cells.Item($u.Row, 5).Value = $u.E
    }
}
